# Update the "as_of_utc" timestamp (column AA) for data rows 2-26 on the
# "Главные" and "Линейные" sheets, reflecting the refreshed publish time.

$wb = $excel.ActiveWorkbook

$oldTimestamp = "2025-12-02 03:04:09"
$newTimestamp = "2025-12-02 07:04:25"

# Sheets 2 ("Главные") and 3 ("Линейные") both carry an as_of_utc column (AA)
# for rows 2 through 26; sheet 1 ("Глоссарий") is a glossary and is untouched.
for ($sheetIndex = 2; $sheetIndex -le 3; $sheetIndex++) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    for ($row = 2; $row -le 26; $row++) {
        $cell = $ws.Cells.Item($row, 27)  # column AA
        if ($cell.Value2 -eq $oldTimestamp) {
            $cell.Value = $newTimestamp
        }
    }
}
